# T4B - Com todas as tabelas
# Recreate the additional "Média / Desvio padrão / u(média)" summary tables
# (diametro, L0 and h/D blocks) on Sheet1, matching the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$center  = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
$vcenter = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter
$thin    = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$xlLeft   = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeLeft
$xlRight  = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight
$xlTop    = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop
$xlBottom = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom

Write-Host "start"

# ---------------------------------------------------------------------------
# Column widths (new columns I and O get explicit widths, like columns
# B, C, F, K:L already had)
# ---------------------------------------------------------------------------
$ws.Columns("I").ColumnWidth = 12.88671875
$ws.Columns("O").ColumnWidth = 13.33203125
